$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row based on the worksheet's used range,
# falling back to the known data extent (rows 2-295) if unavailable.
$lastRow = $ws.UsedRange.Rows.Count
if (-not $lastRow -or $lastRow -lt 295) {
    $lastRow = 295
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46075) {
        $cell.Value2 = 46076
    }
}
